$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wasser1")

# Update the formula in G2 (fix: square the second term inside the SQRT)
$ws.Range("G2").Formula = "=SQRT((1/(B2-D2)-B2/((B2-D2)^2))^2*C2^2+(E2*B2/((B2-D2)^2))^2)"

# Update the formulas in G3:G11 (same correction), one row at a time
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("G$r").Formula = "=SQRT((1/(B$r-D$r)-B$r/((B$r-D$r)^2))^2*C$r^2+(E$r*B$r/((B$r-D$r)^2))^2)"
}

# Update the active cell selection on the sheet
$ws.Range("G8").Select() | Out-Null

$wb.Save()
